# Apply the new table style to the comparison table on slide 16
# (graphicFrame "Google Shape;213;p29", the 3rd shape on that slide).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$sh = $s.Shapes.Item(3)
$tbl = $sh.Table
$tbl.ApplyStyle("{53AA04DF-7E3B-44E5-B1DE-E2153604B7CE}")
